$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where only the "Absent" (H) column flips to 1
$absentRows = @(3, 7, 8, 11, 14, 15, 16, 17, 18)
foreach ($r in $absentRows) {
    $ws.Cells.Item($r, 8).Value = 1   # column H = Absent
}

# Rows where "Total Attendance Count" (D) and "Real" (E) columns flip to 1
$presentRows = @(4, 5, 6, 9, 10, 12, 13)
foreach ($r in $presentRows) {
    $ws.Cells.Item($r, 4).Value = 1   # column D = Total Attendance Count
    $ws.Cells.Item($r, 5).Value = 1   # column E = Real
}
